$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

# Copy formatting (style) from the cell above in column A (date style) so the
# new date cell reuses the existing date number format instead of creating a
# brand-new style entry.
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 42625.885787037034
$ws.Cells.Item($row, 2).Value = 14
$ws.Cells.Item($row, 3).Value = 68
$ws.Cells.Item($row, 4).Value = 31
$ws.Cells.Item($row, 5).Value = 77
$ws.Cells.Item($row, 6).Value = 22
$ws.Cells.Item($row, 7).Value = 11347
$ws.Cells.Item($row, 8).Value = 6199
$ws.Cells.Item($row, 9).Value = 878
$ws.Cells.Item($row, 10).Value = 184
$ws.Cells.Item($row, 11).Value = 85
$ws.Cells.Item($row, 12).Value = 7
$ws.Cells.Item($row, 13).Value = 2
$ws.Cells.Item($row, 14).Value = "Noun"
